$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 13.62512651192116
$ws.Range("C2").Value = 14.01869158878505
$ws.Range("D2").Value = 14.41015089163237
$ws.Range("E2").Value = 13.58103904184166

$ws.Range("B3").Value = 5.966137459807074
$ws.Range("C3").Value = 7.184048480109471
$ws.Range("D3").Value = 2.547654436882318

$ws.Range("B4").Value = 48.21894005212858
$ws.Range("C4").Value = 47.99981415230219
$ws.Range("D4").Value = 47.28610407328006
$ws.Range("E4").Value = 45.65205927332173

$ws.Range("E5").Value = 49.76771196283391

$ws.Range("B6").Value = 43.42524594706942

$ws.Range("E7").Value = 55.54603894698204

$ws.Range("B8").Value = 2.527362816151457
$ws.Range("C8").Value = 3.379581225804629
$ws.Range("D8").Value = 3.621518485479235
$ws.Range("E8").Value = 4.730001771536243
